$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Test Case" / "Test Steps" / "Expected Results" content (columns D, E, G)
# together with each row's height were reshuffled among rows 12-20:
#   new row 12 <- old row 18
#   new row 13 <- old row 19
#   new row 14 <- old row 20
#   new row 15 <- old row 12
#   new row 16 <- old row 13
#   new row 17 <- old row 14
#   new row 18 <- old row 15
#   new row 19 <- old row 16
#   new row 20 <- old row 17
# (columns B, C, F, H, I and the rest of the sheet are untouched.)

# 1) Snapshot the current ("old") D/E/G values and row heights for rows 12-20
#    before any writes happen, so reads never see already-overwritten data.
$oldD = @{}
$oldE = @{}
$oldG = @{}
$oldHt = @{}
for ($r = 12; $r -le 20; $r++) {
    $oldD[$r] = $ws.Range("D$r").Value2
    $oldE[$r] = $ws.Range("E$r").Value2
    $oldG[$r] = $ws.Range("G$r").Value2
    $oldHt[$r] = $ws.Rows.Item($r).RowHeight
}

# Keep a formatting donor for column G's cell style (quote-prefixed "표준 2"
# style) on a scratch cell far outside the used range; writing new text into
# a G12:G20 cell resets its style, so we restore it afterwards from here.
$ws.Range("G12").Copy()
$ws.Range("ZZ1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 2) Define the new-row -> old-row source mapping.
$map = @{
    12 = 18
    13 = 19
    14 = 20
    15 = 12
    16 = 13
    17 = 14
    18 = 15
    19 = 16
    20 = 17
}

# 3) Apply the snapshot values according to the mapping.
foreach ($newRow in 12..20) {
    $srcRow = $map[$newRow]
    $ws.Range("D$newRow").Value = $oldD[$srcRow]
    $ws.Range("E$newRow").Value = $oldE[$srcRow]
    $ws.Range("G$newRow").Value = $oldG[$srcRow]
    $ws.Rows.Item($newRow).RowHeight = $oldHt[$srcRow]
}

# 4) Restore the column-G cell style that got clobbered by the value writes.
$ws.Range("ZZ1").Copy()
$ws.Range("G12:G20").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 5) Clean up the scratch cell used as a formatting donor.
$ws.Range("ZZ1").Clear()
